# "Doing Updates for Financials"
#
# NDBKY yearly financials refresh: the historical figures in columns D:J
# (most-recent-year .. oldest-year) were restated with updated source data.
# A handful of cells that previously held a trailing-year numeric value now
# report "NA" (shared string already used elsewhere on the sheet) because
# that data point is no longer available for the oldest year shown.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NDBKY")

function Set-RowValues {
    # Writes D..J on $row from $values (7 entries). A $null entry leaves the
    # existing cell value untouched (used for the handful of cells in rows
    # 57/58 that were already "NA" and stay that way).
    param($ws, $row, $values)
    $cols = @("D", "E", "F", "G", "H", "I", "J")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        if ($null -ne $values[$i]) {
            $ws.Range("$($cols[$i])$row").Value = $values[$i]
        }
    }
}

# Income Statement
Set-RowValues $ws 8   @(5161000, 5030500, 4132200, 3606500, 3158800, 3065800, 2939000)    # Total Revenue
Set-RowValues $ws 15  @(-156800, -145200, -124100, -114100, -104000, -160900, -95800)     # Others
Set-RowValues $ws 17  @(3494100, 3531400, 2823400, 2341600, 2085800, 2073300, 2068300)    # Total Operating Expenses
Set-RowValues $ws 18  @(1666900, 1499100, 1308800, 1264900, 1073000, 992500, 870700)      # Operating Income or Loss
Set-RowValues $ws 20  @(-535400, -497500, -302600, -328900, -252400, -262400, -275900)    # Total Other Income/Expenses Net
Set-RowValues $ws 21  @(1288300, 1146800, 1130300, 1050100, 924700, 830800, "NA")         # Earnings Before Interest And Taxes
Set-RowValues $ws 23  @(1131500, 1001600, 1006200, 936000, 820600, 730100, 594700)        # Income Before Tax
Set-RowValues $ws 24  @(288500, 271100, 241200, 237700, 206700, 196400, 149000)           # Income Tax Expense
Set-RowValues $ws 26  @(843000, 730600, 765000, 698300, 613900, 533700, 445700)           # Income After Tax
Set-RowValues $ws 27  @(796500, 694400, 734800, 671400, 592000, 510600, 424300)           # Net Income From Continuing Ops
Set-RowValues $ws 32  @(535400, 497500, 302600, 328900, 252400, 262400, 275900)           # Other Items
Set-RowValues $ws 33  @(796500, 694400, 734800, 671400, 592000, 510600, 424300)           # Net Income
Set-RowValues $ws 35  @(796500, 694400, 734800, 671400, 592000, 510600, 424300)           # Net Income Applicable To Common Shares

# Balance Sheet
Set-RowValues $ws 41  @(3265400, 3868700, 3710200, 3038100, 3581100, 4483400, 1741500)    # Cash And Cash Equivalents
Set-RowValues $ws 42  @(9502600, 9155700, 8462400, 7384200, 6407400, 6717400, 2466500)    # Short Term Investments
Set-RowValues $ws 47  @(460700, 450100, 656500, 525700, 75500, 70700, 38900)              # Long Term Investments
Set-RowValues $ws 48  @(610100, 616200, 604200, 541700, 482000, 905100, 907300)           # Property Plant and Equipment
Set-RowValues $ws 49  @(780300, 691100, 617500, 588000, 568200, 1085900, 1066100)         # Goodwill
Set-RowValues $ws 52  @(445600, 410100, 362200, 333900, 219900, 221700, 18800)            # Other Assets
Set-RowValues $ws 54  @(67396300, 66211100, 63449300, 55470300, 51377200, 46809900, 44422600) # Total Assets
Set-RowValues $ws 57  @(1009200, 845800, 554100, 547100, 600900, 555900, $null)           # Accounts Payable (J already "NA")
Set-RowValues $ws 58  @(385700, 423000, 628900, 591900, 326400, $null, $null)             # Short/Current Long Term Debt (I,J already "NA")
Set-RowValues $ws 59  @(78400, 70200, 77700, 61100, 69600, 59500, 56600)                  # Other Current Liabilities
Set-RowValues $ws 61  @(3149300, 3146300, 2454200, 1850700, 1953800, 2076600, 2018000)    # Long Term Debt
Set-RowValues $ws 62  @(329000, 316400, 313600, 297600, 211600, 200100, 193600)           # Other Liabilities
Set-RowValues $ws 66  @(61607600, 60883300, 58325600, 50876500, 47222500, 43136100, 41067900) # Total Liabilities
Set-RowValues $ws 72  @(4243000, 3815500, 3527800, 3411700, 3002900, 2543600, 4333500)    # Retained Earnings
Set-RowValues $ws 76  @(5788800, 5327800, 5123600, 4593800, 4154700, 3673800, 3354800)    # Total Stockholder Equity

# Cash Flow Statement
Set-RowValues $ws 81  @(796500, 694400, 734800, 671400, 592000, 510600, 424300)           # Net Income
Set-RowValues $ws 83  @(156800, 145200, 124100, 114100, 104000, 100700, "NA")             # Depreciation
Set-RowValues $ws 89  @(220300, 294000, 305200, 398800, 833200, 613000, 607500)           # Total Cash Flow From Operating Activities
Set-RowValues $ws 91  @(-226100, -263600, -196500, -170500, -119100, -154100, -116300)    # Capital Expenditures
Set-RowValues $ws 94  @(-419400, -205900, 196500, -648000, -297500, -321900, "NA")        # Total Cash Flows From Investing Activities
Set-RowValues $ws 96  @(-439900, -407700, -395200, -340400, -281900, -242700, -219300)    # Dividends Paid
Set-RowValues $ws 100 @(-407500, 242400, 260600, -146100, -54800, -174900, "NA")          # Total Cash Flows From Financing Activities
Set-RowValues $ws 101 @(-7600, 81600, -20600, -3700, -4400, 1200, "NA")                   # Effect Of Exchange Rate Changes
Set-RowValues $ws 102 @(-614300, 412100, 741700, -399100, 476400, 117400, 388200)         # Change In Cash and Cash Equivalents
